$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Remove all existing hyperlinks first so stale relationship targets do not linger
$ws.Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = "2025-09-17 12:36:24"
$ws.Range("B2").Value = "【急募】Pythonによるスクレイピング開発(既存システム改修)"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5394950"
$ws.Range("G2").Value = 320
$ws.Range("H2").Value = "🔥Python ◆開発,スクレイピング"

# Row 3
$ws.Range("A3").Value = "2025-09-17 12:36:24"
$ws.Range("B3").Value = "自社開発のロジシステムをサポート及び開発できる方募集【PHP, Python, VBA etc】"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5389460"
$ws.Range("G3").Value = 305
$ws.Range("H3").Value = "🔥Python ◆開発 ○PHP"

# Row 4
$ws.Range("A4").Value = "2025-09-17 12:36:24"
$ws.Range("B4").Value = "詳細設計及び、Next.js,node.jsによるWEBアプリケーション開発"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5377709"
$ws.Range("G4").Value = 245
$ws.Range("H4").Value = "🔥Next.js ◆開発,Node.js ◇アプリ"

# Row 5
$ws.Range("A5").Value = "2025-09-17 12:36:24"
$ws.Range("B5").Value = "<Next.js、バックエンド開発> ガントチャートアプリの改修製造"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5379158"
$ws.Range("G5").Value = 225
$ws.Range("H5").Value = "🔥Next.js ◆開発 ◇アプリ"

# Row 6
$ws.Range("A6").Value = "2025-09-17 12:36:24"
$ws.Range("B6").Value = "【急募】保育園シフト自動作成ツールの開発依頼"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5395148"
$ws.Range("G6").Value = 128
$ws.Range("H6").Value = "◆ツール,開発"

# Row 7
$ws.Range("A7").Value = "2025-09-17 12:36:24"
$ws.Range("B7").Value = "【急募】ポーカー大会用アプリ開発のプロを探しています!"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "5,000,000 円 ~ / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5395367"
$ws.Range("G7").Value = 100
$ws.Range("H7").Value = "◆開発 ◇アプリ"

# Row 8
$ws.Range("A8").Value = "2025-09-17 12:36:24"
$ws.Range("B8").Value = "Flutterなどハイブリッドアプリによる業務アプリの開発(スマートウォッチ)"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5379176"
$ws.Range("G8").Value = 100
$ws.Range("H8").Value = "◆開発 ◇アプリ"

# Row 9
$ws.Range("A9").Value = "2025-09-17 12:36:24"
$ws.Range("B9").Value = "【スマホ最優先】キャスト向け会員制Webアプリ開発依頼"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5394619"
$ws.Range("G9").Value = 100
$ws.Range("H9").Value = "◆開発 ◇アプリ"

# Row 10
$ws.Range("A10").Value = "2025-09-17 12:36:24"
$ws.Range("B10").Value = "2026年度新入社員研修Javaメイン講師 (4~6月)"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5395010"
$ws.Range("G10").Value = 85
$ws.Range("H10").Value = "★Java"

# Row 11
$ws.Range("A11").Value = "2025-09-17 12:36:24"
$ws.Range("B11").Value = "2026年度新入社員研修Javaメイン講師"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5394910"
$ws.Range("G11").Value = 85
$ws.Range("H11").Value = "★Java"

# Row 12
$ws.Range("A12").Value = "2025-09-17 12:36:24"
$ws.Range("B12").Value = "業務自動化ブログのTech記事ライター"
$ws.Range("C12").Value = "システム開発"
$ws.Range("D12").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E12").Value = "期限情報なし"
$ws.Range("F12").Value = "https://www.lancers.jp/work/detail/5395045"
$ws.Range("G12").Value = 80
$ws.Range("H12").Value = "◆自動化"

# Row 13
$ws.Range("A13").Value = "2025-09-17 12:36:24"
$ws.Range("B13").Value = "IB報酬を得るための高性能EA開発依頼"
$ws.Range("C13").Value = "システム開発"
$ws.Range("D13").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E13").Value = "期限情報なし"
$ws.Range("F13").Value = "https://www.lancers.jp/work/detail/5392235"
$ws.Range("G13").Value = 68
$ws.Range("H13").Value = "◆開発"

# Row 14
$ws.Range("A14").Value = "2025-09-17 12:36:24"
$ws.Range("B14").Value = "wordpressレンダリングを妨げるリソースの除外"
$ws.Range("C14").Value = "システム開発"
$ws.Range("D14").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E14").Value = "期限情報なし"
$ws.Range("F14").Value = "https://www.lancers.jp/work/detail/5016989"
$ws.Range("G14").Value = 33
$ws.Range("H14").Value = "○WordPress"

# Row 15
$ws.Range("A15").Value = "2025-09-17 12:36:24"
$ws.Range("B15").Value = "ZOHO CRMとZOHO キャンペーン CMSの自動配信システム"
$ws.Range("C15").Value = "システム開発"
$ws.Range("D15").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E15").Value = "期限情報なし"
$ws.Range("F15").Value = "https://www.lancers.jp/work/detail/5390965"
$ws.Range("G15").Value = 33
$ws.Range("H15").ClearContents()

# Row 16
$ws.Range("A16").Value = "2025-09-17 12:36:24"
$ws.Range("B16").Value = "Googleスプレッドシート連携型 データ集計システム改修要件"
$ws.Range("C16").Value = "システム開発"
$ws.Range("D16").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E16").Value = "期限情報なし"
$ws.Range("F16").Value = "https://www.lancers.jp/work/detail/5395383"
$ws.Range("G16").Value = 25
$ws.Range("H16").ClearContents()

# Row 17
$ws.Range("A17").Value = "2025-09-17 12:36:24"
$ws.Range("B17").Value = "OR(operations research)にて最適化の仕組みの構築 (リモート)"
$ws.Range("C17").Value = "システム開発"
$ws.Range("D17").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E17").Value = "期限情報なし"
$ws.Range("F17").Value = "https://www.lancers.jp/work/detail/5367840"
$ws.Range("G17").Value = 25
$ws.Range("H17").ClearContents()

# Row 18
$ws.Range("A18").Value = "2025-09-17 12:36:24"
$ws.Range("B18").Value = "OR(operations research)にて最適化の仕組みの構築(社内常駐)"
$ws.Range("C18").Value = "システム開発"
$ws.Range("D18").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E18").Value = "期限情報なし"
$ws.Range("F18").Value = "https://www.lancers.jp/work/detail/5372984"
$ws.Range("G18").Value = 25
$ws.Range("H18").ClearContents()

# Row 19
$ws.Range("A19").Value = "2025-09-17 12:36:24"
$ws.Range("B19").Value = "新規登録(比較的最近に登録)されたGoogleビジネスプロフィールのデータの取得"
$ws.Range("C19").Value = "システム開発"
$ws.Range("D19").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E19").Value = "期限情報なし"
$ws.Range("F19").Value = "https://www.lancers.jp/work/detail/5395243"
$ws.Range("G19").Value = 18
$ws.Range("H19").ClearContents()

# Row 20
$ws.Range("A20").Value = "2025-09-17 12:36:24"
$ws.Range("B20").Value = "【実績重視】LP用離脱防止ポップアップ作成の依頼"
$ws.Range("C20").Value = "システム開発"
$ws.Range("D20").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E20").Value = "期限情報なし"
$ws.Range("F20").Value = "https://www.lancers.jp/work/detail/5394827"
$ws.Range("G20").Value = 13
$ws.Range("H20").ClearContents()

# Row 21
$ws.Range("A21").Value = "2025-09-17 12:36:24"
$ws.Range("B21").Value = "Android kotlin 画像ファイルのアップロード"
$ws.Range("C21").Value = "システム開発"
$ws.Range("D21").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E21").Value = "期限情報なし"
$ws.Range("F21").Value = "https://www.lancers.jp/work/detail/5395033"
$ws.Range("G21").Value = 10
$ws.Range("H21").ClearContents()

# Row 22
$ws.Range("A22").Value = "2025-09-17 12:36:24"
$ws.Range("B22").Value = "EC CUBE4のメール送信およびファイルダウンロードの改修"
$ws.Range("C22").Value = "システム開発"
$ws.Range("D22").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E22").Value = "期限情報なし"
$ws.Range("F22").Value = "https://www.lancers.jp/work/detail/5394781"
$ws.Range("G22").Value = 10
$ws.Range("H22").ClearContents()

# Re-add hyperlinks for F2:F22
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5394950")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5389460")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5377709")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5379158")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5395148")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5395367")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5379176")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5394619")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5395010")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5394910")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5395045")
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5392235")
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5016989")
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.lancers.jp/work/detail/5390965")
$ws.Hyperlinks.Add($ws.Range("F16"), "https://www.lancers.jp/work/detail/5395383")
$ws.Hyperlinks.Add($ws.Range("F17"), "https://www.lancers.jp/work/detail/5367840")
$ws.Hyperlinks.Add($ws.Range("F18"), "https://www.lancers.jp/work/detail/5372984")
$ws.Hyperlinks.Add($ws.Range("F19"), "https://www.lancers.jp/work/detail/5395243")
$ws.Hyperlinks.Add($ws.Range("F20"), "https://www.lancers.jp/work/detail/5394827")
$ws.Hyperlinks.Add($ws.Range("F21"), "https://www.lancers.jp/work/detail/5395033")
$ws.Hyperlinks.Add($ws.Range("F22"), "https://www.lancers.jp/work/detail/5394781")

Write-Output "edit complete"